# Updates the cryptos price/volume table (columns D=Price, E=Volume(1h))
# for rows 2-51 to the refreshed values from the latest scrape.
# NumberFormat is set to Text ("@") before writing numeric-looking price
# strings (e.g. "36.00", "0.999") so Excel keeps them as literal text
# instead of silently coercing them to numbers and dropping formatting
# such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '43.393.47'
$ws.Cells.Item(2, 5).Value = '  +1.49%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.604.44'
$ws.Cells.Item(3, 5).Value = '  +3.23%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '316.78'
$ws.Cells.Item(5, 5).Value = '  +0.04%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '97.69'
$ws.Cells.Item(6, 5).Value = '  +3.33%  '

$ws.Cells.Item(7, 5).Value = '  -0.13%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 5).Value = '  +2.67%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '36.00'
$ws.Cells.Item(10, 5).Value = '  +0.55%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0817'
$ws.Cells.Item(11, 5).Value = '  +0.83%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.56'
$ws.Cells.Item(12, 5).Value = '  +0.16%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '3.001.50'
$ws.Cells.Item(13, 5).Value = '  +3.13%  '

$ws.Cells.Item(14, 5).Value = '  -0.60%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.575.55'
$ws.Cells.Item(15, 5).Value = '  +1.69%  '

$ws.Cells.Item(16, 5).Value = '  +0.91%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.851'
$ws.Cells.Item(17, 5).Value = '  +0.48%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '43.493.49'
$ws.Cells.Item(18, 5).Value = '  +1.50%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.89'
$ws.Cells.Item(19, 5).Value = '  +2.90%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.75'
$ws.Cells.Item(20, 5).Value = '  -1.85%  '

$ws.Cells.Item(21, 5).Value = '  +0.79%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '69.99'
$ws.Cells.Item(22, 5).Value = '  +0.47%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '255.46'
$ws.Cells.Item(23, 5).Value = '  +1.92%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.99'
$ws.Cells.Item(24, 5).Value = '  +1.89%  '

$ws.Cells.Item(25, 5).Value = '  +3.84%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '27.56'

$ws.Cells.Item(27, 5).Value = '  -0.22%  '

$ws.Cells.Item(28, 5).Value = '  +0.85%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '41.26'
$ws.Cells.Item(29, 5).Value = '  +2.37%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '10.35'
$ws.Cells.Item(30, 5).Value = '  +0.82%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '5.93'
$ws.Cells.Item(31, 5).Value = '  -0.68%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '157.73'
$ws.Cells.Item(32, 5).Value = '  +0.86%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.52'
$ws.Cells.Item(33, 5).Value = '  +7.29%  '

$ws.Cells.Item(34, 5).Value = '  +2.92%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0814'
$ws.Cells.Item(35, 5).Value = '  +3.67%  '

$ws.Cells.Item(36, 5).Value = '  +3.59%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '18.86'
$ws.Cells.Item(37, 5).Value = '  -0.57%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.54'
$ws.Cells.Item(38, 5).Value = '  +11.35%  '

$ws.Cells.Item(39, 5).Value = '  +0.67%  '

$ws.Cells.Item(40, 5).Value = '  +0.26%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '23.13'

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '4.06'
$ws.Cells.Item(42, 5).Value = '  +8.12%  '

$ws.Cells.Item(43, 5).Value = '  +0.45%  '

$ws.Cells.Item(44, 5).Value = '  -0.05%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.27'
$ws.Cells.Item(45, 5).Value = '  -0.50%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.021.05'
$ws.Cells.Item(46, 5).Value = '  +0.04%  '

$ws.Cells.Item(47, 5).Value = '  +2.89%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.855.32'
$ws.Cells.Item(48, 5).Value = '  +3.21%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '84.07'
$ws.Cells.Item(49, 5).Value = '  -1.84%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '75.32'
$ws.Cells.Item(50, 5).Value = '  +2.45%  '

$ws.Cells.Item(51, 5).Value = '  +3.08%  '
